$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = 55
$ws.Range("B56").Value = 1
$ws.Range("C56").Value = "2024-06-16 02:31:56"
$ws.Range("D56").Value = 200
$ws.Range("E56").Value = 14

$ws.Range("A57").Value = 56
$ws.Range("B57").Value = 2
$ws.Range("C57").Value = "2024-06-16 02:31:56"
$ws.Range("D57").Value = 200
$ws.Range("E57").Value = 1
